$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 75
$ws.Range("I2").Value = 212
$ws.Range("J2").Value = 859
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 254
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 180
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 95
$ws.Range("T2").Value = 157
$ws.Range("U2").Value = 15
$ws.Range("V2").Value = 1364
$ws.Range("X2").Value = 1369
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 11
